# Generate Report for Handback
# Updates timestamp / status strings that are shared across rows 2 and 3
# on the "Overview", "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-29 22:17:47"
$wsOverview.Range("G3").Value = "2016-08-29 22:17:47"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-29 22:17:42"
$wsZhCn.Range("H3").Value = "2016-08-29 22:17:42"
$wsZhCn.Range("K2").Value = "2016-08-29 22:17:59"
$wsZhCn.Range("K3").Value = "2016-08-29 22:17:59"

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-29 22:18:12"
$wsDeDe.Range("K3").Value = "2016-08-29 22:18:12"
